# Add 20 new parts rows (rows 8-27) to the blueprint parts list on Sheet1.
# All new rows share material "RENSHAPE" and CNC "Yes", matching the
# existing data pattern (columns B "qty" and E "description" stay blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Item = 7;  Part = "97195A419" },
    @{ Item = 8;  Part = "2489 TR_1" },
    @{ Item = 9;  Part = "2489 CABLE_2" },
    @{ Item = 10; Part = "2489 CABLE_1" },
    @{ Item = 11; Part = "2489 1W_08" },
    @{ Item = 12; Part = "2489 1W_06" },
    @{ Item = 13; Part = "2489 1W_05" },
    @{ Item = 14; Part = "2489 1W_04" },
    @{ Item = 15; Part = "2489 1W_3" },
    @{ Item = 16; Part = "2489 1W_03" },
    @{ Item = 17; Part = "2489 1W_2B_M" },
    @{ Item = 18; Part = "2489 1W_2B" },
    @{ Item = 19; Part = "2489 1W_2A_M" },
    @{ Item = 20; Part = "2489 1W_2A" },
    @{ Item = 21; Part = "2489 1W_02" },
    @{ Item = 22; Part = "2489 1W_02" },
    @{ Item = 23; Part = "2489 1W_02" },
    @{ Item = 24; Part = "2489 1W_1A_M" },
    @{ Item = 25; Part = "2489 1W_1A" },
    @{ Item = 26; Part = "2489 1W_01" }
)

$row = 8
foreach ($data in $newRows) {
    $ws.Range("A$row").Value = $data.Item
    $ws.Range("C$row").Value = $data.Part
    $ws.Range("D$row").Value = "RENSHAPE"
    $ws.Range("F$row").Value = "Yes"
    $row = $row + 1
}

# Leave the selection where the user's cursor ended up after entering data.
$ws.Range("F28").Select()
